$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (moved one month forward: 45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update the three price values in column D
$ws.Range("D29").Value = 348.194
$ws.Range("D30").Value = 368.347
$ws.Range("D31").Value = 396.64
